# Kotak mutual fund portfolio sheet restructuring
# - Inserts 3 new columns (yield, yield to call (ytc), yield to maturity (ytm)) between
#   old "Net Asset Value (NAV)" and "Type" columns.
# - Re-labels / lower-cases several header cells.
# - Moves the "Rating/Industry" data (old column C) into the new "industry" column (D),
#   clearing the old column C (now "coupon", left blank for this data set).
# - Moves the old "Type" (H), "Scheme" (I) and "AmcName" (J) data into the new
#   K / L / M columns, clearing the old H/I/J (now blank "yield"/"yield to call
#   (ytc)"/"yield to maturity (ytm)" columns), and rewrites the Type value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Create the three brand-new header cells (I1, J1 used to hold "Scheme" /
#    "AmcName" - those get pushed out to L1/M1). Copy the existing bold/
#    bordered header style (from A1) onto the new K1/L1/M1 cells before
#    populating them so the style index matches the rest of row 1.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 9))
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 11))
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 12))
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 13))

# ---------------------------------------------------------------------------
# 2. Header row text (columns A-M).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value  = "name of instrument"
$ws.Cells.Item(1, 2).Value  = "isin"
$ws.Cells.Item(1, 3).Value  = "coupon"
$ws.Cells.Item(1, 4).Value  = "industry"
$ws.Cells.Item(1, 5).Value  = "quantity"
$ws.Cells.Item(1, 6).Value  = "market value (mkt)"
$ws.Cells.Item(1, 7).Value  = "% to net assets (nav)"
$ws.Cells.Item(1, 8).Value  = "yield"
$ws.Cells.Item(1, 9).Value  = "yield to call (ytc)"
$ws.Cells.Item(1, 10).Value = "yield to maturity (ytm)"
$ws.Cells.Item(1, 11).Value = "Type"
$ws.Cells.Item(1, 12).Value = "Scheme"
$ws.Cells.Item(1, 13).Value = "AmcName"

# ---------------------------------------------------------------------------
# 3. Data rows (2-44): shuffle the per-row values into their new homes.
# ---------------------------------------------------------------------------
$newType = "Equity & Equity related   NAN nan nan nan nan nan"

for ($r = 2; $r -le 44; $r++) {

    # --- read old values before overwriting anything (use Value2 - the
    #     plain `.Value` getter is not usable for reads in this runtime) ---
    $industry = $ws.Cells.Item($r, 3).Value2    # old C : Rating/Industry
    $scheme   = $ws.Cells.Item($r, 9).Value2    # old I : Scheme (TIF/TCH)
    $amcName  = $ws.Cells.Item($r, 10).Value2   # old J : AmcName

    # --- column C ("coupon") is now blank for this data set ---
    $ws.Cells.Item($r, 3).Value = ""

    # --- column D ("industry") takes over the old Rating/Industry text ---
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $industry
    $ws.Cells.Item($r, 4).Style = "Normal"

    # --- columns H/I/J ("yield" / "yield to call (ytc)" / "yield to maturity (ytm)") are blank ---
    $ws.Cells.Item($r, 8).Value  = ""
    $ws.Cells.Item($r, 9).Value  = ""
    $ws.Cells.Item($r, 10).Value = ""

    # --- column K ("Type") gets the rewritten type text ---
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = $newType
    $ws.Cells.Item($r, 11).Style = "Normal"

    # --- column L ("Scheme") takes over the old Scheme text ---
    $ws.Cells.Item($r, 12).NumberFormat = "@"
    $ws.Cells.Item($r, 12).Value = $scheme
    $ws.Cells.Item($r, 12).Style = "Normal"

    # --- column M ("AmcName") takes over the old AmcName text ---
    $ws.Cells.Item($r, 13).NumberFormat = "@"
    $ws.Cells.Item($r, 13).Value = $amcName
    $ws.Cells.Item($r, 13).Style = "Normal"
}
